$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("index")

# "Table 1".."Table 6" -> "Table_1".."Table_6" (spaces not supported by bookdown)
for ($i = 1; $i -le 6; $i++) {
    $ws.Range("B" + (18 + $i)).Value = "Table_" + $i
}

$ws.Range("B24").Select()
